$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header, formatted like the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New column values (save/era data)
$values = @(1, 1, 1, 1, 0, 1, 1, 1, 0, 1, 0, 1, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
